$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: mapping of after-row -> target values, derived from the row
# that originally (pre-edit) held that date, reflecting a re-sort of the
# weekly Fruta/Hortaliza records (rows 2-15) by Fecha.

# Row 2 <- was row 7
$ws.Cells.Item(2, 4).Value = 44589
$ws.Cells.Item(2, 13).Value = 60
$ws.Cells.Item(2, 14).Value = 6000
$ws.Cells.Item(2, 15).Value = 6000
$ws.Cells.Item(2, 16).Value = 6000
$ws.Cells.Item(2, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(2, 19).Value = 3000

# Row 3 <- was row 9
$ws.Cells.Item(3, 4).Value = 44587
$ws.Cells.Item(3, 13).Value = 165
$ws.Cells.Item(3, 14).Value = 6500
$ws.Cells.Item(3, 15).Value = 7000
$ws.Cells.Item(3, 16).Value = 6742
$ws.Cells.Item(3, 18).Value = "Provincia de Linares"
$ws.Cells.Item(3, 19).Value = 3371

# Row 4 <- was row 15
$ws.Cells.Item(4, 4).Value = 44211
$ws.Cells.Item(4, 13).Value = 45
$ws.Cells.Item(4, 14).Value = 6000
$ws.Cells.Item(4, 15).Value = 6000
$ws.Cells.Item(4, 16).Value = 6000
$ws.Cells.Item(4, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(4, 19).Value = 3000

# Row 5 <- was row 3
$ws.Cells.Item(5, 4).Value = 44585
$ws.Cells.Item(5, 13).Value = 160
$ws.Cells.Item(5, 14).Value = 6500
$ws.Cells.Item(5, 15).Value = 7000
$ws.Cells.Item(5, 16).Value = 6750
$ws.Cells.Item(5, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(5, 19).Value = 3375

# Row 6 <- was row 8
$ws.Cells.Item(6, 4).Value = 44209
$ws.Cells.Item(6, 13).Value = 58
$ws.Cells.Item(6, 14).Value = 6000
$ws.Cells.Item(6, 15).Value = 6000
$ws.Cells.Item(6, 16).Value = 6000
$ws.Cells.Item(6, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(6, 19).Value = 3000

# Row 7 <- was row 10
$ws.Cells.Item(7, 4).Value = 44214
$ws.Cells.Item(7, 13).Value = 48
$ws.Cells.Item(7, 14).Value = 6000
$ws.Cells.Item(7, 15).Value = 6000
$ws.Cells.Item(7, 16).Value = 6000
$ws.Cells.Item(7, 18).Value = "Provincia de Linares"
$ws.Cells.Item(7, 19).Value = 3000

# Row 8 <- was row 12
$ws.Cells.Item(8, 4).Value = 44586
$ws.Cells.Item(8, 13).Value = 80
$ws.Cells.Item(8, 14).Value = 7000
$ws.Cells.Item(8, 15).Value = 7000
$ws.Cells.Item(8, 16).Value = 7000
$ws.Cells.Item(8, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(8, 19).Value = 3500

# Row 9 <- was row 2
$ws.Cells.Item(9, 4).Value = 44588
$ws.Cells.Item(9, 13).Value = 160
$ws.Cells.Item(9, 14).Value = 6500
$ws.Cells.Item(9, 15).Value = 7000
$ws.Cells.Item(9, 16).Value = 6750
$ws.Cells.Item(9, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(9, 19).Value = 3375

# Row 10 <- was row 4
$ws.Cells.Item(10, 4).Value = 44606
$ws.Cells.Item(10, 13).Value = 45
$ws.Cells.Item(10, 14).Value = 7000
$ws.Cells.Item(10, 15).Value = 7000
$ws.Cells.Item(10, 16).Value = 7000
$ws.Cells.Item(10, 18).Value = "Provincia de Linares"
$ws.Cells.Item(10, 19).Value = 3500

# Row 11 <- was row 6
$ws.Cells.Item(11, 4).Value = 44627
$ws.Cells.Item(11, 13).Value = 45
$ws.Cells.Item(11, 14).Value = 6000
$ws.Cells.Item(11, 15).Value = 6000
$ws.Cells.Item(11, 16).Value = 6000
$ws.Cells.Item(11, 18).Value = "Provincia de Linares"
$ws.Cells.Item(11, 19).Value = 3000

# Row 12 <- was row 13
$ws.Cells.Item(12, 4).Value = 44582
$ws.Cells.Item(12, 13).Value = 150
$ws.Cells.Item(12, 14).Value = 6000
$ws.Cells.Item(12, 15).Value = 6500
$ws.Cells.Item(12, 16).Value = 6233
$ws.Cells.Item(12, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(12, 19).Value = 3116

# Row 13 <- was row 11
$ws.Cells.Item(13, 4).Value = 44628
$ws.Cells.Item(13, 13).Value = 40
$ws.Cells.Item(13, 14).Value = 6000
$ws.Cells.Item(13, 15).Value = 6000
$ws.Cells.Item(13, 16).Value = 6000
$ws.Cells.Item(13, 18).Value = "Provincia de Linares"
$ws.Cells.Item(13, 19).Value = 3000

# Row 14 <- was row 5
$ws.Cells.Item(14, 4).Value = 44614
$ws.Cells.Item(14, 13).Value = 45
$ws.Cells.Item(14, 14).Value = 6000
$ws.Cells.Item(14, 15).Value = 6000
$ws.Cells.Item(14, 16).Value = 6000
$ws.Cells.Item(14, 18).Value = "Provincia de Linares"
$ws.Cells.Item(14, 19).Value = 3000

# Row 15 <- was row 14
$ws.Cells.Item(15, 4).Value = 44592
$ws.Cells.Item(15, 13).Value = 30
$ws.Cells.Item(15, 14).Value = 8000
$ws.Cells.Item(15, 15).Value = 8000
$ws.Cells.Item(15, 16).Value = 8000
$ws.Cells.Item(15, 18).Value = "Provincia de Linares"
$ws.Cells.Item(15, 19).Value = 4000
